$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark that currently wraps the picture run in
#    the image paragraph (bookmarkStart before the drawing run,
#    bookmarkEnd after it).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# 2. Append " Maecenas sit amet consequat diam. ..." to the paragraph
#    that currently ends with "...makes it hard to remember."
$p4 = $d.Paragraphs.Item(4)
$insertPos = $p4.Range.End - 1
$r = $d.Range($insertPos, $insertPos)
$r.InsertAfter(" ")

$insertPos2 = $insertPos + 1
$r2 = $d.Range($insertPos2, $insertPos2)
$r2.InsertAfter("Maecenas sit amet consequat diam. Phasellus molestie lorem eros, elementum consectetur purus placerat in. Fusce sed tristique risus, sit amet pellentesque erat. Donec a urna sed nisi aliquet posuere vel a nulla. Integer consequat hendrerit vulputate. Curabitur et tellus nec metus tempus tempor. In id sodales turpis. Nulla lacus est, convallis sit amet ligula non, ornare placerat tellus. Duis posuere lacus sit amet tempus scelerisque. Etiam id cursus neque. Morbi vitae ligula massa. Phasellus lacinia felis hendrerit sem commodo, id tempus ipsum imperdiet. Suspendisse eu ante mi. Sed molestie, ipsum at feugiat venenatis, magna lectus tristique ex, non semper mauris tortor sit amet massa. Aliquam eu erat turpis. Nam dignissim velit in tellus interdum, vel feugiat felis elementum.")

# 3. Split off a brand-new empty paragraph right after that paragraph
#    (inserted just before its own paragraph mark, so no stray run is
#    left behind) and re-home the _GoBack bookmark there, leaving the
#    original trailing empty paragraph untouched as the final paragraph.
$p4b = $d.Paragraphs.Item(4)
$splitPos = $p4b.Range.End - 1
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertBefore([char]13)

$newPara = $d.Paragraphs.Item(5)
$d.Bookmarks.Add("_GoBack", $newPara.Range)
